$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Matryoshka" row (row 13); rows below shift up by one.
$ws.Rows.Item(13).Delete()

# After the delete, row 15 is "Kendama" (the last existing data row) and rows
# 16 onward are already-empty placeholder rows, so we just fill them in
# directly (no row insertion needed - that would also push down the sheet's
# trailing placeholder rows, which the authored edit does not do).

# Data for the new "National Flags Hat" entries.
$flags = @(
    @("America", "America.png", "America_climb.png"),
    @("England", "England.png", "England_climb.png"),
    @("Germany", "Germany.png", "Germany_climb.png"),
    @("Greece", "Greece.png", "Greece_climb.png"),
    @("Italy", "Italy.png", "Italy_climb.png"),
    @("Jamaica", "Jamaica.png", "Jamaica_climb.png"),
    @("Japan", "Japan.png", "Japan_climb.png"),
    @("Pakistan", "Pakistan.png", "Pakistan_climb.png"),
    @("Panama", "Panama.png", "Panama_climb.png"),
    @("Russia", "Russia.png", "Russia_climb.png"),
    @("Seychelles", "Seychelles.png", "Seychelles_climb.png"),
    @("Spain", "Spain.png", "Spain_climb.png"),
    @("Sri Lanka", "SriLanka.png", "SriLanka_climb.png"),
    @("Turkey", "Turkey.png", "Turkey_climb.png"),
    @("Ukraine", "Ukraine.png", "Ukraine_climb.png")
)

$row = 16
foreach ($flag in $flags) {
    $ws.Cells.Item($row, 1).Value2 = $flag[0]
    $ws.Cells.Item($row, 2).Value2 = "9on"
    $ws.Cells.Item($row, 3).Value2 = "nationalFlagHat"
    $ws.Cells.Item($row, 4).Value2 = $flag[1]
    $ws.Cells.Item($row, 6).Value2 = $flag[2]
    $row = $row + 1
}

# Update the active selection to match the authored state.
$ws.Range("F31").Select()
